$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.959.40"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = "'1.879.30"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.86%  '
$ws.Range('D4').Value = "'0.9989"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'242.61"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.32%  '
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').Value = "'0.4975"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.33%  '
$ws.Range('D8').Value = "'0.2925"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('D9').Value = "'0.06632"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.62%  '
$ws.Range('D10').Value = "'1.879.31"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('D11').Value = "'16.74"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.01%  '
$ws.Range('D12').Value = "'0.07240"
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = "'0.6680"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.09%  '
$ws.Range('D14').Value = "'86.35"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.57%  '
$ws.Range('D15').Value = "'4.890"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.42%  '
$ws.Range('D16').Value = "'29.949.95"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('D17').Value = "'0.000007903"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.01%  '
$ws.Range('D18').Value = "'0.9981"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('E19').Value = '  -1.96%  '
$ws.Range('D20').Value = "'2.120.34"
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Value = "'0.9974"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').Value = "'4.772"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.76%  '
$ws.Range('D23').Value = "'5.656"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.46%  '
$ws.Range('D24').Value = "'9.062"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('D25').Value = "'149.32"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.99%  '
$ws.Range('D26').Value = "'141.96"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.35%  '
$ws.Range('D27').Value = "'17.16"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.76%  '
$ws.Range('D28').Value = "'1.912"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.36%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = "'4.180"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.97%  '
$ws.Range('D31').Value = "'0.08780"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.61%  '
$ws.Range('E32').Value = '  -2.02%  '
$ws.Range('D33').Value = "'0.05067"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.00%  '
$ws.Range('D34').Value = "'0.7091"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.18%  '
$ws.Range('D35').Value = "'1.107"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.63%  '
$ws.Range('D36').Value = "'2.665"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.65%  '
$ws.Range('D37').Value = "'0.01775"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.49%  '
$ws.Range('D38').Value = "'2.686"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.69%  '
$ws.Range('D39').Value = "'2.177"
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Value = "'0.9308"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.63%  '
$ws.Range('D41').Value = "'5.800"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.61%  '
$ws.Range('E42').Value = '  -1.59%  '
$ws.Range('D43').Value = "'0.9977"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = "'102.32"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.87%  '
$ws.Range('D45').Value = "'7.482"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.25%  '
$ws.Range('D46').Value = "'0.1259"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.94%  '
$ws.Range('E47').Value = '  -1.88%  '
$ws.Range('E48').Value = '  -3.25%  '
$ws.Range('D49').Value = "'0.3766"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.85%  '
$ws.Range('D50').Value = "'8.235"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.53%  '
$ws.Range('D51').Value = "'55.90"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.79%  '
